$wb = $excel.ActiveWorkbook

# --- Sheet 1: 台指期換倉成本計算 ---
$ws1 = $wb.Worksheets.Item("台指期換倉成本計算")
$ws1.Rows.Item(2).Insert()
$ws1.Range("A2").Value = "日期：2021/12/23"
# Force B2 to be stored as plain text (not auto-converted to a number)
$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = "202202"
$ws1.Range("B2").ClearFormats()
$ws1.Range("C2").Value = 17922
$ws1.Range("D2").Value = 2735
$ws1.Range("E2").Value = 4283358
$ws1.Range("F2").Value = 17663

# --- Sheet 2: 散戶多空力道 ---
$ws2 = $wb.Worksheets.Item("散戶多空力道")
$ws2.Rows.Item(2).Insert()
$ws2.Range("A2").Value = "日期：2021/12/23"
$ws2.Range("B2").Value = -0.03

# --- Sheet 3: 三大法人買賣金額 ---
$ws3 = $wb.Worksheets.Item("三大法人買賣金額")
$ws3.Rows.Item(2).Insert()
$ws3.Range("A2").Value = "110年12月23日"
$ws3.Range("B2").Value = 180.23
$ws3.Range("C2").Value = -60.4

# --- Sheet 4: 大盤多空點位 ---
$ws4 = $wb.Worksheets.Item("大盤多空點位")
$ws4.Rows.Item(2).Insert()
$ws4.Range("A2").Value = "110年12月23日"
$ws4.Range("B2").Value = 17931.41

# --- Sheet 5: 期貨大額交易人未沖銷部位 ---
$ws5 = $wb.Worksheets.Item("期貨大額交易人未沖銷部位")
$ws5.Rows.Item(2).Insert()
# Force A2 to be stored as plain text (not auto-converted to a date serial)
$ws5.Range("A2").NumberFormat = "@"
$ws5.Range("A2").Value = "2021/12/23"
$ws5.Range("A2").ClearFormats()
$ws5.Range("B2").Value = 45936
$ws5.Range("C2").Value = 52465
$ws5.Range("D2").Value = -387
$ws5.Range("E2").Value = -584
$ws5.Range("F2").Value = 22282
$ws5.Range("G2").Value = 46148
$ws5.Range("H2").Value = 1085
$ws5.Range("I2").Value = 757
$ws5.Range("J2").Value = -23866
$ws5.Range("K2").Value = 328
$ws5.Range("L2").Value = -1472
$ws5.Range("M2").Value = -1341
$ws5.Range("N2").Value = -131
